# RESERVE_settings.xlsx update
#
# Commit: "Created persistence forecasting method for solar using
#          theoretical clear sky output."
#
# Changes:
#   1. "RESERVE Settings" sheet: add a new LATITUDE parameter row (with its
#      Description and Value) right before the existing LONGITUDE row, so a
#      solar clear-sky model has the latitude it needs.
#   2. "Input Data Settings" sheet: the Solar/Forecast row's Data Source
#      switches from the old 15-min-ahead forecast CSV to the new
#      "persistence" method (the same method already used for the Wind
#      forecast row).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("RESERVE Settings")
$ws2 = $wb.Worksheets.Item("Input Data Settings")

# --- "RESERVE Settings": insert LATITUDE row above LONGITUDE (row 6) ---
$ws1.Rows("6:6").Insert()
$ws1.Range("A6").Value = "LATITUDE"
$ws1.Range("B6").Value = "Approximate latitude of power system"
$ws1.Range("C6").Value = 36.6777

# New page setup orientation recorded for this sheet
$ws1.PageSetup.Orientation = 1

# --- "Input Data Settings": Solar Forecast row now sources from persistence ---
$ws2.Range("A7").Value = "persistence"

# --- Restore on-screen selections (sheet2 stays the active/visible tab) ---
[void]$ws1.Range("B14").Select()
[void]$ws2.Range("A11").Select()
